# Append the March 4th attendance record to the sheet as a new row.
# The Date/Entry-Time values must stay literal text (not get silently
# converted into Excel date/time serials), so force text formatting
# before assigning them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-03-04"
$ws.Range("B2").Value = "Vasanth Kumar"
$ws.Range("C2").Value = "19:06:39"
$ws.Range("D2").Value = ""

[void]$ws.Range("A1").Select()
